# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N) for specific
# leve rows across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets with the
# latest marketboard figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2620.2
$ws.Range("I40").Value = 1024.5
$ws.Range("J40").Value = 3684
$ws.Range("K40").Value = 1024.5
$ws.Range("L40").Value = 3684
$ws.Range("M40").Value = -849.5
$ws.Range("N40").Value = -4034

$ws.Range("H64").Value = 24323.459
$ws.Range("I64").Value = 3526.8823
$ws.Range("J64").Value = 74829.42999999999
$ws.Range("K64").Value = 3526.8823
$ws.Range("L64").Value = 74829.42999999999
$ws.Range("M64").Value = -3278.8823
$ws.Range("N64").Value = -75325.42999999999

$ws.Range("H67").Value = 24323.459
$ws.Range("I67").Value = 3526.8823
$ws.Range("J67").Value = 74829.42999999999
$ws.Range("K67").Value = 3526.8823
$ws.Range("L67").Value = 74829.42999999999
$ws.Range("M67").Value = -2668.8823
$ws.Range("N67").Value = -76545.42999999999

$ws.Range("H74").Value = 3350.158
$ws.Range("I74").Value = 3352.9443
$ws.Range("J74").Value = 3300
$ws.Range("K74").Value = 3352.9443
$ws.Range("L74").Value = 3300
$ws.Range("M74").Value = -2416.9443
$ws.Range("N74").Value = -5172

$ws.Range("H76").Value = 27030258
$ws.Range("I76").Value = 31253190
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 31253190
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -31252875
$ws.Range("N76").Value = -4130

$ws.Range("H77").Value = 3350.158
$ws.Range("I77").Value = 3352.9443
$ws.Range("J77").Value = 3300
$ws.Range("K77").Value = 16764.7215
$ws.Range("L77").Value = 16500
$ws.Range("M77").Value = -12084.7215
$ws.Range("N77").Value = -25860

$ws.Range("H79").Value = 27030258
$ws.Range("I79").Value = 31253190
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 31253190
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -31252098
$ws.Range("N79").Value = -5684

$ws.Range("H129").Value = 897.07245
$ws.Range("I129").Value = 394
$ws.Range("J129").Value = 963.0492
$ws.Range("K129").Value = 1182
$ws.Range("L129").Value = 2889.1476
$ws.Range("M129").Value = 3818
$ws.Range("N129").Value = -12889.1476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2350.84
$ws.Range("I61").Value = 2132.1904
$ws.Range("J61").Value = 3498.75
$ws.Range("K61").Value = 2132.1904
$ws.Range("L61").Value = 3498.75
$ws.Range("M61").Value = -1920.1904
$ws.Range("N61").Value = -3922.75

$ws.Range("H63").Value = 2354.5557
$ws.Range("I63").Value = 2378.9333
$ws.Range("J63").Value = 2232.6667
$ws.Range("K63").Value = 2378.9333
$ws.Range("L63").Value = 2232.6667
$ws.Range("M63").Value = -1692.9333
$ws.Range("N63").Value = -3604.6667

$ws.Range("H66").Value = 2354.5557
$ws.Range("I66").Value = 2378.9333
$ws.Range("J66").Value = 2232.6667
$ws.Range("K66").Value = 11894.6665
$ws.Range("L66").Value = 11163.3335
$ws.Range("M66").Value = -8462.666500000001
$ws.Range("N66").Value = -18027.3335

$ws.Range("H74").Value = 1341.825
$ws.Range("I74").Value = 1234.4062
$ws.Range("J74").Value = 1771.5
$ws.Range("K74").Value = 1234.4062
$ws.Range("L74").Value = 1771.5
$ws.Range("M74").Value = -360.4061999999999

$ws.Range("H77").Value = 1341.825
$ws.Range("I77").Value = 1234.4062
$ws.Range("J77").Value = 1771.5
$ws.Range("K77").Value = 6172.030999999999
$ws.Range("L77").Value = 8857.5
$ws.Range("M77").Value = -1804.030999999999

$ws.Range("H136").Value = 2350.84
$ws.Range("I136").Value = 2132.1904
$ws.Range("J136").Value = 3498.75
$ws.Range("K136").Value = 6396.5712
$ws.Range("L136").Value = 10496.25
$ws.Range("M136").Value = -3846.5712
$ws.Range("N136").Value = -15596.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1594522
$ws.Range("I105").Value = 1992452.5
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 1992452.5
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -1990705.5
$ws.Range("N105").Value = -6294

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2269.1
$ws.Range("I103").Value = 1108.3334
$ws.Range("J103").Value = 2766.5715
$ws.Range("K103").Value = 3325.0002
$ws.Range("L103").Value = 8299.7145
$ws.Range("M103").Value = -2446.0002
$ws.Range("N103").Value = -10057.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6066.6665
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6066.6665
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6066.6665
$ws.Range("N70").Value = -6606.6665
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 6066.6665
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6066.6665
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6066.6665
$ws.Range("N73").Value = -7938.6665
$ws.Range("M73").ClearContents()

$ws.Range("H80").Value = 9143.799999999999
$ws.Range("I80").Value = 2917.5
$ws.Range("J80").Value = 13294.667
$ws.Range("K80").Value = 2917.5
$ws.Range("L80").Value = 13294.667
$ws.Range("M80").Value = -1919.5
$ws.Range("N80").Value = -15290.667

$ws.Range("H83").Value = 9143.799999999999
$ws.Range("I83").Value = 2917.5
$ws.Range("J83").Value = 13294.667
$ws.Range("K83").Value = 14587.5
$ws.Range("L83").Value = 66473.33499999999
$ws.Range("M83").Value = -9595.5
$ws.Range("N83").Value = -76457.33499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2304.6206
$ws.Range("I7").Value = 2432.95
$ws.Range("J7").Value = 2019.4445
$ws.Range("K7").Value = 2432.95
$ws.Range("L7").Value = 2019.4445
$ws.Range("M7").Value = -2320.95
$ws.Range("N7").Value = -2243.4445

$ws.Range("H126").Value = 2304.6206
$ws.Range("I126").Value = 2432.95
$ws.Range("J126").Value = 2019.4445
$ws.Range("K126").Value = 7298.849999999999
$ws.Range("L126").Value = 6058.333500000001
$ws.Range("M126").Value = -4828.849999999999
$ws.Range("N126").Value = -10998.3335

$ws.Range("H136").Value = 557531.8
$ws.Range("I136").Value = 835022.3
$ws.Range("J136").Value = 2550.8333
$ws.Range("K136").Value = 2505066.9
$ws.Range("L136").Value = 7652.499899999999
$ws.Range("M136").Value = -2502516.9
$ws.Range("N136").Value = -12752.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1896056.5
$ws.Range("I136").Value = 3109114.2
$ws.Range("J136").Value = 501040.06
$ws.Range("K136").Value = 9327342.600000001
$ws.Range("L136").Value = 1503120.18
$ws.Range("M136").Value = -9324792.600000001
$ws.Range("N136").Value = -1508220.18

$ws.Range("H139").Value = 50897.895
$ws.Range("I139").Value = 29800
$ws.Range("J139").Value = 52070
$ws.Range("K139").Value = 29800
$ws.Range("L139").Value = 52070
$ws.Range("M139").Value = -24660
$ws.Range("N139").Value = -62350
